$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Insert a new row at row 5 ("LP solver (linprog or gurobi)" / "gurobi"),
# pushing the existing rows (previously 5-14) down to 6-15.
$ws.Range("A5:B5").Insert()
$ws.Range("A5").Value = "LP solver (linprog or gurobi)"
$ws.Range("B5").Value = "gurobi"

# Match the author's final selection state recorded in the saved file.
$ws.Range("A5:B5").Select()
